$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.589.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.17%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.618.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.69%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.90%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.14%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.622.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.17%  "

# Row 10
$ws.Range("E10").Value = "  -2.32%  "

# Row 11
$ws.Range("E11").Value = "  -0.22%  "

# Row 12
$ws.Range("E12").Value = "  -1.26%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.127"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.95%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.076.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.61%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.583.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.17%  "

# Row 17
$ws.Range("E17").Value = "  -1.32%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.618.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "345.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.29%  "

# Row 22
$ws.Range("E22").Value = "  -0.61%  "

# Row 23
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "

# Row 25
$ws.Range("E25").Value = "  -1.79%  "

# Row 26
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.163"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.38%  "

# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.33%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0796"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.84%  "

# Row 30
$ws.Range("E30").Value = "  +0.10%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.73%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.77%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.96%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.77%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.974"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.87%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.22%  "

# Row 37
$ws.Range("E37").Value = "  -0.41%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.98%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.832"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.10%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.46%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "278.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.47%  "

# Row 43
$ws.Range("E43").Value = "  +0.22%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0978"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.90%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.91%  "

# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.596"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.29%  "

# Row 47
$ws.Range("E47").Value = "  +0.58%  "

# Row 48
$ws.Range("E48").Value = "  -3.53%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0228"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.56%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.975.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.47%  "
